$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that is bumped by one
# day (45837 -> 45838) for every data row (rows 2 through 43).
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45837) {
        $cell.Value2 = 45838
    }
}
